$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the InfiniBand-vs-Ethernet "speeds" paragraph with the new HPC
#    Advisory Council study language, and color it dark red (C00000) to
#    match the author's edit.
# ---------------------------------------------------------------------------
$oldText = "InfiniBand is considerably faster that Ethernet.  The typical ISP delivers an Ethernet solution that is less than 1Gb/sec.  Often much less.  InfiniBand can deliver (btw, at a steep price) speeds as fast 40Gb/sec.  That" + [char]8217 + "s more than an order of magnitude faster than a typical Ethernet solution."
$newText = "InfiniBand is considerably faster that Ethernet.  A comprehensive study published by the High Performance Computing Advisory Council produced concrete metrics that demonstrated that Infiniband delivers 600% better low-latency  and 370% better throughput performance  than Ethernet (10GE). "

$findRng = $d.Content
$found = $findRng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Color the whole rewritten paragraph (text + paragraph mark) dark red.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*comprehensive study published by the High Performance Computing*") {
        $p.Range.Font.Color = 192
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Add the new "HPC Advisory Council Study ..." citation paragraph (plus a
#    trailing blank paragraph) right after the existing RDMA/Wikipedia
#    hyperlink paragraph, citing the study referenced above.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Remote_Direct_Memory_Access*") {
        $anchorIndex = $i
        break
    }
}

$anchorRng = $d.Paragraphs.Item($anchorIndex).Range
$anchorRng.Collapse(0)
$anchorRng.InsertParagraphAfter()

$newParaIndex = $anchorIndex + 1
$newRng = $d.Paragraphs.Item($newParaIndex).Range
$newRng.Collapse(1)
$citationText = "HPC Advisory Council Study Infiniband vs. Ethernet 10GE. http://www.hpcadvisorycouncil.com/pdf/IB_and_10GigE_in_HPC.pdf"
$newRng.InsertAfter($citationText)

# Wrap the URL portion in a real hyperlink.
$urlAddress = "http://www.hpcadvisorycouncil.com/pdf/IB_and_10GigE_in_HPC.pdf"
$urlRng = $d.Paragraphs.Item($newParaIndex).Range
$urlFound = $urlRng.Find.Execute($urlAddress, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($urlRng, $urlAddress, "", "", $urlAddress) | Out-Null

# Match the surrounding small (8pt / sz=16) font used throughout this block.
$sizeRng = $d.Paragraphs.Item($newParaIndex).Range
$sizeRng.Find.Execute($urlAddress, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sizeRng.Font.Size = 8

# Blank paragraph separating the citation from the pre-existing blank line.
$blankAnchor = $d.Paragraphs.Item($newParaIndex).Range
$blankAnchor.Collapse(0)
$blankAnchor.InsertParagraphAfter()

Write-Output "done"
